$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 57: same four-column shape row 56 currently has (before its
#     own restyle below), so clone that formatting first.
$ws.Range("A56:D56").Copy() | Out-Null
$ws.Range("A57:D57").PasteSpecial(-4122) | Out-Null
$ws.Rows(57).RowHeight = 23.25

# --- Row 56: A56's key ("cwl_warn_processor") switches to the non-wrapping
#     "warn/log id" look already used by cells like A55 (style shared with
#     A36..A47 etc: gold Cascadia Code, vertical-center, no wrap). Grab that
#     formatting via copy/paste so we land on the identical style bucket.
$ws.Range("A55").Copy() | Out-Null
$ws.Range("A56").PasteSpecial(-4122) | Out-Null

# Column A: loc key name
$ws.Range("A57").Value = "cwl_log_ele_gain"

# Column B stays blank (already cleared by the format-only paste above)

# Column C: English source string
$ws.Range("C57").Value = "auto gained ability id: {0} on {1}"

# Column D: ZH-CN translation, built from four differently-fonted runs
#   "已添加能力: " (宋体) + "{0} " (Cascadia Code) + "至" (宋体) + " {1}" (Cascadia Code)
$ws.Range("D57").Value = "已添加能力: {0} 至 {1}"
$ws.Range("D57").Characters(1, 7).Font.Name = "宋体"
$ws.Range("D57").Characters(1, 7).Font.ColorIndex = -4105
$ws.Range("D57").Characters(8, 4).Font.Name = "Cascadia Code"
$ws.Range("D57").Characters(8, 4).Font.ColorIndex = -4105
$ws.Range("D57").Characters(12, 1).Font.Name = "宋体"
$ws.Range("D57").Characters(12, 1).Font.ColorIndex = -4105
$ws.Range("D57").Characters(13, 4).Font.Name = "Cascadia Code"
$ws.Range("D57").Characters(13, 4).Font.ColorIndex = -4105

$excel.CutCopyMode = $false
$ws.Range("D57").Select() | Out-Null
